$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 14 content (order matches shared-string insertion order: A, B, D, C, E, F)
$ws.Range("A14").Value = "Checking to see if all of the values inside of a column are unique (for the purpose of validating primary keys)"
$ws.Range("B14").Value = "Run a script that checks each value and counts how many identical values are inside of the column, the script then returns the count of each value that it found"
$ws.Range("D14").Value = "SELECT PCID, COUNT(PCID) `nFROM PCALLOCATION`nGROUP BY PCID`nHAVING COUNT(PCID) > 1;"
$ws.Range("C14").Value = "For the purpose of primary keys, the expected result of this test would be for every value to be unique(no value having COUNT > 1). For the purpose of this test, the PCALLOCATION table will be used for the PCID column"
$ws.Range("E14").Value = "No value having COUNT > 1"
$ws.Range("F14").Value = "Checks to see that every value used for a column is unique"

# Match style of row 13: A/B/D/F use wrap text (s="1"), C/E use center + wrap (s="4")
$ws.Range("A14").WrapText = $true
$ws.Range("B14").WrapText = $true
$ws.Range("D14").WrapText = $true
$ws.Range("F14").WrapText = $true

$ws.Range("C14").WrapText = $true
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("E14").WrapText = $true
$ws.Range("E14").HorizontalAlignment = -4108

$ws.Rows.Item(14).RowHeight = 85.5

$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D22").Select()
